{"js": "const body = context.document.body;\n\nasync function replaceOnce(searchText, replacementText) {\n  const results = body.search(searchText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error('Text not found: ' + searchText);\n  }\n  results.items[0].insertText(replacementText, Word.InsertLocation.replace);\n  await context.sync();\n}\n\nawait replaceOnce(\"\u00bfQu\u00e9 es Azure Application Gateway?\", \"Informaci\u00f3n general del producto\");\nawait replaceOnce(\"Contoso CipherGuard Sentinel X7 es un producto de seguridad avanzado y resistente dise\u00f1ado meticulosamente para fortalecer la infraestructura de la red inform\u00e1tica frente a una amplia gama de amenazas y vulnerabilidades.\", \"Contoso CipherGuard Sentinel X7 es un producto de seguridad avanzado y resistente dise\u00f1ado cuidadosamente para fortalecer la infraestructura de red inform\u00e1tica frente a un espectro de amenazas y vulnerabilidades.\");\nawait replaceOnce(\" el uso de un firewall de inspecci\u00f3n con estado, Contoso CipherGuard Sentinel X7 emplea t\u00e9cnicas de inspecci\u00f3n profunda de paquetes.\", \" al usar un firewall de inspecci\u00f3n con estado, Contoso CipherGuard Sentinel X7 emplea t\u00e9cnicas de inspecci\u00f3n profunda de paquetes.\");\nawait replaceOnce(\"con tecnolog\u00eda de algoritmos de aprendizaje autom\u00e1tico, nuestro IDPS supervisa continuamente los patrones y anomal\u00edas del tr\u00e1fico de red.\", \"con tecnolog\u00eda de algoritmos de aprendizaje autom\u00e1tico, nuestro IDPS supervisa continuamente las anomal\u00edas y los patrones de tr\u00e1fico de red.\");\nawait replaceOnce(\"Contoso CipherGuard Sentinel X7 admite protocolos VPN est\u00e1ndar del sector, como IPsec y OpenVPN.\", \"Contoso CipherGuard Sentinel X7 admite los protocolos VPN est\u00e1ndar del sector, como IPsec y OpenVPN.\");\nawait replaceOnce(\"Seguridad de los puntos de conexi\u00f3n:\", \"Seguridad de puntos de conexi\u00f3n:\");\nawait replaceOnce(\" el uso de un enfoque de defensa multicapa, nuestro m\u00f3dulo de seguridad de puntos de conexi\u00f3n incorpora funcionalidades antivirus, antimalware y prevenci\u00f3n de intrusiones basadas en host.\", \" al usar un enfoque de defensa multicapa, nuestro m\u00f3dulo de seguridad de puntos de conexi\u00f3n incorpora funcionalidades antivirus, antimalware y de prevenci\u00f3n de intrusiones basadas en host.\");\nawait replaceOnce(\"Autenticaci\u00f3n de usuario y control\", \"Autenticaci\u00f3n de usuario y control de acceso\");\nawait replaceOnce(\" de acceso: Contoso CipherGuard Sentinel X7 admite mecanismos de autenticaci\u00f3n multifactor (MFA), incluida la autenticaci\u00f3n biom\u00e9trica y la integraci\u00f3n de tarjetas inteligentes.\", \": Contoso CipherGuard Sentinel X7 admite mecanismos de autenticaci\u00f3n multifactor (MFA), incluida la autenticaci\u00f3n biom\u00e9trica y la integraci\u00f3n de tarjetas inteligentes.\");\nawait replaceOnce(\" Cuatro n\u00facleos de 2,5 GHz o superior con compatibilidad con aceleraci\u00f3n de hardware\", \" cuatro n\u00facleos de 2,5 GHz o superior con compatibilidad con aceleraci\u00f3n de hardware\");\nawait replaceOnce(\" 16 GB como m\u00ednimo, ECC (c\u00f3digo de correcci\u00f3n de errores) recomendado\", \" m\u00ednimo 16 GB, se recomienda ECC (c\u00f3digo de correcci\u00f3n de errores)\");\nawait replaceOnce(\" 200 GB como m\u00ednimo, SSD para un rendimiento \u00f3ptimo\", \" m\u00ednimo 200 GB, SSD para un rendimiento \u00f3ptimo\");\nawait replaceOnce(\" actualizaciones automatizadas para fuentes de inteligencia sobre amenazas y revisiones de seguridad normales\", \" actualizaciones automatizadas para fuentes de inteligencia sobre amenazas y parches peri\u00f3dicos de seguridad\");\nawait replaceOnce(\" interoperabilidad con Cisco, Juniper y otros principales proveedores de redes\", \" interoperabilidad con Cisco, Juniper y otros proveedores principales de redes\");\nawait replaceOnce(\" realice una evaluaci\u00f3n completa de vulnerabilidades de red, incluidas las pruebas de penetraci\u00f3n y el an\u00e1lisis de riesgos.\", \" realizar una evaluaci\u00f3n completa de las vulnerabilidades de red, incluidas las pruebas de penetraci\u00f3n y el an\u00e1lisis de riesgos.\");\nawait replaceOnce(\" implemente Contoso CipherGuard Sentinel X7 en servidores dedicados o m\u00e1quinas virtuales, lo que garantiza un uso \u00f3ptimo del hardware y la asignaci\u00f3n de recursos.\", \" implementar Contoso CipherGuard Sentinel X7 en servidores dedicados o m\u00e1quinas virtuales, lo que garantiza un uso \u00f3ptimo del hardware y de la asignaci\u00f3n de recursos.\");\nawait replaceOnce(\" personalice las directivas de seguridad, los controles de acceso y las reglas de firewall en funci\u00f3n de los requisitos de la organizaci\u00f3n.\", \" personalizar las directivas de seguridad, los controles de acceso y las reglas de firewall en funci\u00f3n de los requisitos de la organizaci\u00f3n.\");\nawait replaceOnce(\" ejecute un plan de pruebas exhaustivo, incluidos escenarios de ataque simulados y pruebas de carga, para validar la eficacia y el rendimiento de la soluci\u00f3n.\", \" ejecutar un plan de pruebas exhaustivo, incluyendo escenarios de ataque simulados y pruebas de carga, para validar la eficacia y el rendimiento de la soluci\u00f3n.\");\nawait replaceOnce(\" proporcione sesiones de aprendizaje detalladas para el personal de TI, que cubre operaciones diarias, procedimientos de respuesta a incidentes y tareas de mantenimiento.\", \" proporcionar sesiones de aprendizaje detalladas para el personal de TI, que cubran las operaciones diarias, los procedimientos de respuesta a incidentes y las tareas de mantenimiento.\");\nawait replaceOnce(\" Contoso garantiza actualizaciones continuas en el producto, incorporando la inteligencia sobre amenazas y las mejoras de seguridad m\u00e1s recientes.\", \" Contoso garantiza actualizaciones continuas del producto, incorporando la inteligencia sobre amenazas y las mejoras de seguridad m\u00e1s recientes.\");\nawait replaceOnce(\" Contoso proporciona un equipo de soporte t\u00e9cnico dedicado de 24/7 para garantizar la asistencia r\u00e1pida de cualquier problema t\u00e9cnico o consulta relacionado con Contoso CipherGuard Sentinel X7 .\", \" Contoso proporciona un equipo de soporte t\u00e9cnico dedicado las 24 horas del d\u00eda, los 7 d\u00edas de la semana, para garantizar la asistencia r\u00e1pida para cualquier consulta o problema t\u00e9cnico relacionado con Contoso CipherGuard Sentinel X7 .\");\n", "ps1": "$d = $word.ActiveDocument\n\nfunction Replace-Text($findText, $replaceText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n    # wdFindContinue=1, wdReplaceOne=1\n    $find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 1) | Out-Null\n}\n\nReplace-Text \"\u00bfQu\u00e9 es Azure Application Gateway?\" \"Informaci\u00f3n general del producto\"\nReplace-Text \"Contoso CipherGuard Sentinel X7 es un producto de seguridad avanzado y resistente dise\u00f1ado meticulosamente para fortalecer la infraestructura de la red inform\u00e1tica frente a una amplia gama de amenazas y vulnerabilidades.\" \"Contoso CipherGuard Sentinel X7 es un producto de seguridad avanzado y resistente dise\u00f1ado cuidadosamente para fortalecer la infraestructura de red inform\u00e1tica frente a un espectro de amenazas y vulnerabilidades.\"\nReplace-Text \" el uso de un firewall de inspecci\u00f3n con estado, Contoso CipherGuard Sentinel X7 emplea t\u00e9cnicas de inspecci\u00f3n profunda de paquetes.\" \" al usar un firewall de inspecci\u00f3n con estado, Contoso CipherGuard Sentinel X7 emplea t\u00e9cnicas de inspecci\u00f3n profunda de paquetes.\"\nReplace-Text \"con tecnolog\u00eda de algoritmos de aprendizaje autom\u00e1tico, nuestro IDPS supervisa continuamente los patrones y anomal\u00edas del tr\u00e1fico de red.\" \"con tecnolog\u00eda de algoritmos de aprendizaje autom\u00e1tico, nuestro IDPS supervisa continuamente las anomal\u00edas y los patrones de tr\u00e1fico de red.\"\nReplace-Text \"Contoso CipherGuard Sentinel X7 admite protocolos VPN est\u00e1ndar del sector, como IPsec y OpenVPN.\" \"Contoso CipherGuard Sentinel X7 admite los protocolos VPN est\u00e1ndar del sector, como IPsec y OpenVPN.\"\nReplace-Text \"Seguridad de los puntos de conexi\u00f3n:\" \"Seguridad de puntos de conexi\u00f3n:\"\nReplace-Text \" el uso de un enfoque de defensa multicapa, nuestro m\u00f3dulo de seguridad de puntos de conexi\u00f3n incorpora funcionalidades antivirus, antimalware y prevenci\u00f3n de intrusiones basadas en host.\" \" al usar un enfoque de defensa multicapa, nuestro m\u00f3dulo de seguridad de puntos de conexi\u00f3n incorpora funcionalidades antivirus, antimalware y de prevenci\u00f3n de intrusiones basadas en host.\"\nReplace-Text \"Autenticaci\u00f3n de usuario y control\" \"Autenticaci\u00f3n de usuario y control de acceso\"\nReplace-Text \" de acceso: Contoso CipherGuard Sentinel X7 admite mecanismos de autenticaci\u00f3n multifactor (MFA), incluida la autenticaci\u00f3n biom\u00e9trica y la integraci\u00f3n de tarjetas inteligentes.\" \": Contoso CipherGuard Sentinel X7 admite mecanismos de autenticaci\u00f3n multifactor (MFA), incluida la autenticaci\u00f3n biom\u00e9trica y la integraci\u00f3n de tarjetas inteligentes.\"\nReplace-Text \" Cuatro n\u00facleos de 2,5 GHz o superior con compatibilidad con aceleraci\u00f3n de hardware\" \" cuatro n\u00facleos de 2,5 GHz o superior con compatibilidad con aceleraci\u00f3n de hardware\"\nReplace-Text \" 16 GB como m\u00ednimo, ECC (c\u00f3digo de correcci\u00f3n de errores) recomendado\" \" m\u00ednimo 16 GB, se recomienda ECC (c\u00f3digo de correcci\u00f3n de errores)\"\nReplace-Text \" 200 GB como m\u00ednimo, SSD para un rendimiento \u00f3ptimo\" \" m\u00ednimo 200 GB, SSD para un rendimiento \u00f3ptimo\"\nReplace-Text \" actualizaciones automatizadas para fuentes de inteligencia sobre amenazas y revisiones de seguridad normales\" \" actualizaciones automatizadas para fuentes de inteligencia sobre amenazas y parches peri\u00f3dicos de seguridad\"\nReplace-Text \" interoperabilidad con Cisco, Juniper y otros principales proveedores de redes\" \" interoperabilidad con Cisco, Juniper y otros proveedores principales de redes\"\nReplace-Text \" realice una evaluaci\u00f3n completa de vulnerabilidades de red, incluidas las pruebas de penetraci\u00f3n y el an\u00e1lisis de riesgos.\" \" realizar una evaluaci\u00f3n completa de las vulnerabilidades de red, incluidas las pruebas de penetraci\u00f3n y el an\u00e1lisis de riesgos.\"\nReplace-Text \" implemente Contoso CipherGuard Sentinel X7 en servidores dedicados o m\u00e1quinas virtuales, lo que garantiza un uso \u00f3ptimo del hardware y la asignaci\u00f3n de recursos.\" \" implementar Contoso CipherGuard Sentinel X7 en servidores dedicados o m\u00e1quinas virtuales, lo que garantiza un uso \u00f3ptimo del hardware y de la asignaci\u00f3n de recursos.\"\nReplace-Text \" personalice las directivas de seguridad, los controles de acceso y las reglas de firewall en funci\u00f3n de los requisitos de la organizaci\u00f3n.\" \" personalizar las directivas de seguridad, los controles de acceso y las reglas de firewall en funci\u00f3n de los requisitos de la organizaci\u00f3n.\"\nReplace-Text \" ejecute un plan de pruebas exhaustivo, incluidos escenarios de ataque simulados y pruebas de carga, para validar la eficacia y el rendimiento de la soluci\u00f3n.\" \" ejecutar un plan de pruebas exhaustivo, incluyendo escenarios de ataque simulados y pruebas de carga, para validar la eficacia y el rendimiento de la soluci\u00f3n.\"\nReplace-Text \" proporcione sesiones de aprendizaje detalladas para el personal de TI, que cubre operaciones diarias, procedimientos de respuesta a incidentes y tareas de mantenimiento.\" \" proporcionar sesiones de aprendizaje detalladas para el personal de TI, que cubran las operaciones diarias, los procedimientos de respuesta a incidentes y las tareas de mantenimiento.\"\nReplace-Text \" Contoso garantiza actualizaciones continuas en el producto, incorporando la inteligencia sobre amenazas y las mejoras de seguridad m\u00e1s recientes.\" \" Contoso garantiza actualizaciones continuas del producto, incorporando la inteligencia sobre amenazas y las mejoras de seguridad m\u00e1s recientes.\"\nReplace-Text \" Contoso proporciona un equipo de soporte t\u00e9cnico dedicado de 24/7 para garantizar la asistencia r\u00e1pida de cualquier problema t\u00e9cnico o consulta relacionado con Contoso CipherGuard Sentinel X7 .\" \" Contoso proporciona un equipo de soporte t\u00e9cnico dedicado las 24 horas del d\u00eda, los 7 d\u00edas de la semana, para garantizar la asistencia r\u00e1pida para cualquier consulta o problema t\u00e9cnico relacionado con Contoso CipherGuard Sentinel X7 .\"\n"}
